# Resumen_Semanal.xlsx weekly update
# - "Super Linea" sheet: refresh per-category sales vs. goal figures for the week
#   (category order reshuffled + new amounts / % cumplimiento + traffic-light fill)
# - "Semana Sucursal" sheet: header rename (Trafico/Conversion -> metas) +
#   add the new "TIENDA ONLINE" branch row
# - sharedStrings gains "TIENDA ONLINE" as a brand-new destinatario/branch

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Super Linea" sheet — reorder categories & refresh D/E/F for rows 2-9
# ---------------------------------------------------------------------------
$wsLinea = $wb.Worksheets.Item("Super Linea")

$lineaCats  = @("LINEA BLANCA","OTRAS LINEAS","HOGAR","COMPUTO","CELULARES","ELECTRONICA","SERVICIOS","AUTOMOTRIZ")
$lineaD     = @(93873.5414023716, 37634.854169395, 38144.7088221312, 28482.8795219057, 39012.6866394289, 73699.7265201087, 144895.8732323615, 17349.8496081647)
$lineaE     = @(76972.38, 21605.58, 43018.02, 23069.79, 66409.71, 66796.86, 146745.51, 17620.11)
$lineaF     = @(121.9574, 174.1904, 88.6714, 123.4639, 58.7454, 110.3341, 98.7395, 98.4661)
# 1 = verde (cumple), 2 = amarillo (alerta), 3 = rojo (no cumple) -- matches the
# existing conditional fills already baked into the sheet's style table
$lineaStyle = @(1, 1, 2, 1, 3, 1, 2, 2)

# Stable reference cells (outside the rows we are about to touch) that already
# carry each of the three traffic-light styles, so PasteSpecial formats reuses
# the existing style index instead of minting duplicate style entries.
$styleRefs = @{
    1 = $wb.Worksheets.Item("Semana Sucursal").Range("K2")
    2 = $wb.Worksheets.Item("Semana Sucursal").Range("J2")
    3 = $wb.Worksheets.Item("Semana Sucursal").Range("G4")
}

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $wsLinea.Range("C$r").Value = $lineaCats[$i]
    $wsLinea.Range("D$r").Value = $lineaD[$i]
    $wsLinea.Range("E$r").Value = $lineaE[$i]
    $wsLinea.Range("F$r").Value = $lineaF[$i]

    $styleRefs[$lineaStyle[$i]].Copy()
    $wsLinea.Range("F$r").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) "Semana Sucursal" sheet — rename headers + append TIENDA ONLINE row
# ---------------------------------------------------------------------------
$wsSuc = $wb.Worksheets.Item("Semana Sucursal")

$wsSuc.Range("H1").Value = "Meta Trafico"
$wsSuc.Range("I1").Value = "Meta Conversion (%)"

$wsSuc.Range("A7").Value = 6
$wsSuc.Range("B7").Value = 2022
$wsSuc.Range("C7").Value = 42
$wsSuc.Range("D7").Value = "TIENDA ONLINE"
$wsSuc.Range("E7").Value = 3229.125321
$wsSuc.Range("F7").Value = 3312.96
$wsSuc.Range("G7").Value = 97.4694
$wsSuc.Range("H7").Value = 5125
$wsSuc.Range("I7").Value = 0.0041
$wsSuc.Range("J7").Value = 0
$wsSuc.Range("K7").Value = 0
$wsSuc.Range("L7").Value = 23
$wsSuc.Range("M7").Value = 0

# G7 = amarillo (style 2), J7/K7 = rojo (style 3)
$wsSuc.Range("J2").Copy()
$wsSuc.Range("G7").PasteSpecial(-4122)
$wsSuc.Range("G4").Copy()
$wsSuc.Range("J7").PasteSpecial(-4122)
$wsSuc.Range("G4").Copy()
$wsSuc.Range("K7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
